$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.976.00'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.234.62'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.11'
$ws.Range('E5').Value = '  -4.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.51'
$ws.Range('E6').Value = '  -5.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.568'
$ws.Range('E7').Value = '  -0.72%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -4.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.68'
$ws.Range('E10').Value = '  -5.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  -2.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.18'
$ws.Range('E12').Value = '  -4.55%  '
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('D14').Value = '2.574.18'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '2.235.78'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.822'
$ws.Range('E16').Value = '  -3.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.59'
$ws.Range('E17').Value = '  -5.19%  '
$ws.Range('D18').Value = '43.859.60'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').Value = '0.0₃0956'
$ws.Range('E19').Value = '  -1.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.06'
$ws.Range('E20').Value = '  -10.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.23'
$ws.Range('E21').Value = '  -3.12%  '
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.98'
$ws.Range('E24').Value = '  -5.63%  '
$ws.Range('E25').Value = '  -5.19%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.97'
$ws.Range('E27').Value = '  -6.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.45'
$ws.Range('E28').Value = '  -3.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.16'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('E30').Value = '  -2.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '152.88'
$ws.Range('E32').Value = '  -4.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0799'
$ws.Range('E33').Value = '  -4.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.20'
$ws.Range('E34').Value = '  +4.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.56'
$ws.Range('E35').Value = '  -4.11%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('E37').Value = '  -6.81%  '
$ws.Range('E38').Value = '  -9.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.03'
$ws.Range('E39').Value = '  -8.32%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.80'
$ws.Range('E40').Value = '  -8.62%  '
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.33'
$ws.Range('E41').Value = '  -8.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0300'
$ws.Range('E42').Value = '  -4.29%  '
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').Value = '1.727.70'
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '85.46'
$ws.Range('E45').Value = '  +5.80%  '
$ws.Range('E46').Value = '  -4.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '99.56'
$ws.Range('E47').Value = '  -3.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.89'
$ws.Range('E48').Value = '  -5.06%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.04'
$ws.Range('E49').Value = '  -2.87%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.46'
$ws.Range('E50').Value = '  -8.10%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.12'
$ws.Range('E51').Value = '  -5.33%  '

$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'

Write-Output "Applied 95 cell updates (30 forced to text)"
